$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; Col=4; Value="245.56"},
    @{Row=2; Col=5; Value="-0.50%"},
    @{Row=3; Col=4; Value="28.47"},
    @{Row=3; Col=5; Value="-3.34%"},
    @{Row=4; Col=4; Value="5.285"},
    @{Row=4; Col=5; Value="1.62%"},
    @{Row=5; Col=4; Value="0.05708"},
    @{Row=5; Col=5; Value="-0.46%"},
    @{Row=6; Col=4; Value="6.638"},
    @{Row=6; Col=5; Value="1.08%"},
    @{Row=7; Col=4; Value="3.211"},
    @{Row=7; Col=5; Value="3.54%"},
    @{Row=9; Col=4; Value="0.8945"},
    @{Row=9; Col=5; Value="3.19%"},
    @{Row=10; Col=4; Value="0.1396"},
    @{Row=10; Col=5; Value="2.26%"},
    @{Row=11; Col=4; Value="0.07095"},
    @{Row=11; Col=5; Value="0.07%"},
    @{Row=12; Col=4; Value="0.03162"},
    @{Row=12; Col=5; Value="4.93%"},
    @{Row=13; Col=4; Value="0.09227"},
    @{Row=13; Col=5; Value="-1.72%"},
    @{Row=14; Col=4; Value="0.001534"},
    @{Row=14; Col=5; Value="-0.28%"},
    @{Row=15; Col=4; Value="0.0005968"},
    @{Row=15; Col=5; Value="-0.26%"},
    @{Row=16; Col=4; Value="0.005930"},
    @{Row=16; Col=5; Value="-1.19%"},
    @{Row=17; Col=4; Value="3.495"},
    @{Row=17; Col=5; Value="0.02%"},
    @{Row=18; Col=4; Value="2.172"},
    @{Row=18; Col=5; Value="-0.51%"},
    @{Row=19; Col=4; Value="0.3166"},
    @{Row=19; Col=5; Value="-0.58%"},
    @{Row=20; Col=4; Value="0.03340"},
    @{Row=20; Col=5; Value="0.56%"},
    @{Row=21; Col=4; Value="0.1306"},
    @{Row=21; Col=5; Value="1.23%"},
    @{Row=22; Col=4; Value="3.488"},
    @{Row=22; Col=5; Value="0.07%"},
    @{Row=23; Col=4; Value="0.04064"},
    @{Row=23; Col=5; Value="-1.89%"},
    @{Row=24; Col=4; Value="0.1379"},
    @{Row=24; Col=5; Value="-0.11%"},
    @{Row=25; Col=4; Value="0.001223"},
    @{Row=25; Col=5; Value="-0.26%"},
    @{Row=26; Col=5; Value="-16.87%"},
    @{Row=27; Col=5; Value="-0.88%"},
    @{Row=40; Col=4; Value="0.03786"},
    @{Row=40; Col=5; Value="0.84%"},
    @{Row=41; Col=5; Value="-0.52%"},
    @{Row=42; Col=5; Value="-35.60%"},
    @{Row=43; Col=4; Value="0.002420"},
    @{Row=43; Col=5; Value="-1.28%"},
    @{Row=44; Col=4; Value="0.009446"},
    @{Row=44; Col=5; Value="-0.57%"},
    @{Row=45; Col=4; Value="0.00005279"},
    @{Row=45; Col=5; Value="-0.38%"},
    @{Row=46; Col=5; Value="-0.09%"},
    @{Row=47; Col=4; Value="0.08907"},
    @{Row=47; Col=5; Value="56.17%"},
    @{Row=48; Col=4; Value="0.002259"},
    @{Row=48; Col=5; Value="-0.82%"},
    @{Row=49; Col=4; Value="0.00002099"},
    @{Row=49; Col=5; Value="-0.09%"},
    @{Row=50; Col=4; Value="0.0001999"},
    @{Row=50; Col=5; Value="-0.09%"}
)

foreach ($chg in $changes) {
    $cell = $ws.Cells.Item($chg.Row, $chg.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
}

Write-Output "Applied $($changes.Count) cell updates"
